$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three cells whose text changed. Setting Value first (in the
# order the new shared strings should appear) and then touching Font.Name
# mirrors how Excel materializes an explicit font/style for a freshly
# retyped cell, producing a dedicated cellXf shared by all three cells.
$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Range("B10").Font.Name = "Calibri"

$ws.Range("B9").Value = "0 (312)32-55-46"
$ws.Range("B9").Font.Name = "Calibri"

$ws.Range("B4").Value = "3.2.2 Коэффициент неонатальной смертности"
$ws.Range("B4").Font.Name = "Calibri"

$ws.Range("B8").Select()
